$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 7, 11, 15, 19)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "Logarítmica"
}
